# Append a new scraped listing ("paperspace/comfyui") as the new row 15,
# pushing the previous row 15 ("Stable Diffusion...") down to row 16, and
# refresh every row's "fetched at" timestamp (column A) to the new run
# time - per the commit "Append: 2025-11-16 01:25 JST".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-11-16 01:25:11"

# --- 1. Insert a fresh row above the current last row (row 15). This
#        pushes the existing row 15 ("Stable Diffusion...") -- values,
#        formatting and all -- down to row 16. ---
$ws.Rows.Item(15).Insert()

# --- 2. The insert leaves row 16's hyperlink behind (Excel hyperlinks
#        don't auto-follow a shifted row here), so re-create it pointing
#        at the same URL the row already had. ---
$oldUrl = "https://www.lancers.jp/work/detail/5432055"
$ws.Hyperlinks.Add($ws.Range("F16"), $oldUrl) | Out-Null
$ws.Range("F16").Style = "Hyperlink"

# --- 3. Write the new listing into the now-empty row 15. ---
$newUrl = "https://www.lancers.jp/work/detail/5434935"
$ws.Range("A15").Value = $newTimestamp
$ws.Range("B15").Value = "paperspaceで、comfyuiが動くようにして欲しい。"
$ws.Range("C15").Value = "システム開発"
$ws.Range("D15").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E15").Value = "期限情報なし"
$ws.Range("F15").Value = $newUrl
$ws.Range("G15").Value = 10

# --- 4. Refresh the fetched-at timestamp for every other data row. ---
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
$ws.Cells.Item(16, 1).Value = $newTimestamp
